# Commit: "Changed the suite name to module name."
# The "Test Suite" sheet listed module rows where column A held a
# "Suite" label (e.g. "C Suite", "D Suite") instead of the short module
# name used elsewhere in the column (e.g. "IAM", "Search", "Watchlist").
# Fix the two stray "Suite" labels so column A is consistent with the
# rest of the sheet, matching the module described in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Row 4: "C Suite" -> "Authoring" (column B already says "Authoring module")
$ws.Range("A4").Value = "Authoring"

# Row 5: "D Suite" -> "Profile" (column B already says "Profile module")
$ws.Range("A5").Value = "Profile"

# Move the active selection to A5, matching the cell that was edited last.
$ws.Range("A5").Select()
